$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# Fix the Avg_Time_ms values for the first two result rows (5000 and 10000
# row-count runs) so the sheet reflects the corrected sort timings. The
# scatter chart already plots Data!$D$2:$D$8, so it picks up these
# corrected values automatically - no need to touch the chart itself.
$ws.Range("D2").Value = 1.0506697
$ws.Range("D3").Value = 2.3988462
